# "finish with result file"
#
# Semantic changes made in Excel:
#  1. Sheet "Radar chart": rename the two abbreviated headers to their
#     full form, and correct the "Физическое состояние" data point.
#  2. Sheet "Trend chart": turn the "январь" text label in A1 into a real
#     date (01-Jan-2017) and resize column A to fit.
#  3. Sheet "Name": the pilot number was entered in column B by mistake;
#     move it over to column A (the id column) where the column width was
#     already set up for it.
#
# All the shared-string renumbering visible in the raw XML for the other
# "Trend chart with ..." sheets is a pure side effect of the string table
# shrinking/growing and needs no direct edits on those sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Radar chart
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Radar chart")
$ws1.Range("A1").Value = "Физическое состояние"
$ws1.Range("B1").Value = "Психологическое состояние"
$ws1.Range("B2").Value = 90
$ws1.Range("B4").Select()

# ---------------------------------------------------------------------
# 2. Trend chart
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Trend chart")
$ws2.Range("A1").NumberFormat = "mm-dd-yy"
$ws2.Range("A1").Value = (Get-Date -Year 2017 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0)
$ws2.Columns.Item(1).AutoFit()
$ws2.Range("B6").Select()

# ---------------------------------------------------------------------
# 3. Name
# ---------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item("Name")
$header = $ws7.Range("B1").Value()
$pilotNo = $ws7.Range("B2").Value()
$ws7.Range("A1").Value = $header
$ws7.Range("A2").Value = $pilotNo
$ws7.Range("B1").ClearContents()
$ws7.Range("B2").ClearContents()
$ws7.Range("B2").Select()
